$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer 1 (Pearson Edexcel logo, docPr id="1"): image1.png -> image2.png
$f1 = $sec.Footers.Item(1)
$f1.Range.InlineShapes.Item(1).Name = "image2.png"

# Footer 2 (Pearson Edexcel logo, docPr id="2"): image1.png -> image2.png
$f2 = $sec.Footers.Item(2)
$f2.Range.InlineShapes.Item(1).Name = "image2.png"

# Header 2 (BTec logo, docPr id="3"): image2.jpg -> image1.jpg
$h2 = $sec.Headers.Item(2)
$h2.Range.InlineShapes.Item(1).Name = "image1.jpg"
